$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.22%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.15%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.703"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.79%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06192"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.27%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'0.82%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'-0.94%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9165"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.84%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'1.34%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.04647"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.95%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07085"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.82%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03135"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.87%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09038"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.001544"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.85%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006147"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.50%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006070"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.47%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.458"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.03%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'0.43%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.195"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.22%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-1.01%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'1.60%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.124"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.17%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04229"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.18%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.06%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-5.76%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'0.06%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001601"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-6.48%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.04076"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'6.26%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'-0.21%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'8.03%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-9.73%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-8.30%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005136"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.06%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.06%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'0.1677"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'28.72%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.06%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("E50").Style = "Normal"
